# Generate Report for Archive
# The file "9f10c461-6f73-47ff-a86d-7b8154ea58c7.md" moved from "Ready for
# handoff" back to "In Translation" status. Update the Status column for
# that row (row 4) on the Overview sheet (Status for both zh-cn and de-de
# columns) as well as on the per-locale zh-cn and de-de sheets (Status
# column there).

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B4").Value = "In Translation"
$overview.Range("C4").Value = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B4").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B4").Value = "In Translation"
